$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 792 (shifts existing rows 792.. down by one, dimension grows to D834)
$ws.Rows.Item(792).Insert()

# Populate the newly inserted row with the new data point (2026/02/14 07:01 UTC push).
# Column A holds the date as plain text (matching the rest of the column). Temporarily
# force a text number format so the "yyyy/mm/dd"-shaped string isn't auto-parsed into a
# date serial, then restore the default (unstyled) look so no stray style is left behind.
$ws.Cells.Item(792, 1).NumberFormat = "@"
$ws.Cells.Item(792, 1).Value = "2026/02/14"
$ws.Cells.Item(792, 1).Style = "Normal"

$ws.Cells.Item(792, 2).Value = "土"
$ws.Cells.Item(792, 3).Value = 14
$ws.Cells.Item(792, 4).Value = 22
